$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.162.94"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "2.618.34"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.46"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.10%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "99.35"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.33%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.597"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("E8").Value = "  +0.05%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "38.98"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("E11").Value = "  +2.06%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "54.19"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.25%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "8.10"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.00%  "
$ws.Range("D14").Value = "3.021.56"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "2.615.49"
$ws.Range("E16").Value = "  +1.21%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.917"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.37%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.87"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "46.403.15"
$ws.Range("E19").Value = "  +0.91%  "
$ws.Range("E20").Value = "  +1.94%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.77"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.95%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "12.78"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.87%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "291.57"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +15.81%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "73.21"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.00%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.06"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("E26").Value = "  +3.67%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "29.74"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +8.44%  "
$ws.Range("E28").Value = "  -0.11%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "10.81"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +5.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "39.10"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  -1.42%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.28"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.11%  "
$ws.Range("E34").Value = "  -0.63%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "158.73"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("E36").Value = "  +0.83%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0842"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("E38").Value = "  -4.08%  "
$ws.Range("E39").Value = "  +6.60%  "
$ws.Range("E40").Value = "  +2.18%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "15.74"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.12%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0331"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.33%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.58"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.67%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "21.69"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +10.83%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.04"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "2.129.93"
$ws.Range("E46").Value = "  +4.24%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "97.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +7.10%  "
$ws.Range("E48").Value = "  +0.04%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.50"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.95%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "110.06"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.91%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.201"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.91%  "
